$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.635.15'
$ws.Range('E2').Value = '  +3.70%  '
$ws.Range('D3').Value = '3.499.79'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.90'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.35%  '
$ws.Range('D8').Value = '3.498.00'
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('E9').Value = '  +4.47%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').Value = '  +3.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.439'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.13%  '
$ws.Range('D13').Value = '4.106.14'
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.79%  '
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '66.635.97'
$ws.Range('E17').Value = '  +3.62%  '
$ws.Range('D18').Value = '3.506.33'
$ws.Range('E18').Value = '  +1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.06'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '389.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('E26').Value = '  +3.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.13%  '
$ws.Range('E28').Value = '  +1.83%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +4.10%  '
$ws.Range('E31').Value = '  +5.10%  '
$ws.Range('E32').Value = '  +2.56%  '
$ws.Range('E33').Value = '  +2.34%  '
$ws.Range('E34').Value = '  +4.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +5.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.61'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.24%  '
$ws.Range('E38').Value = '  +2.63%  '
$ws.Range('E39').Value = '  +2.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.71'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0748'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.46'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.824.34'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.38%  '
$ws.Range('E47').Value = '  +2.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '354.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.52%  '
$ws.Range('E49').Value = '  +4.94%  '
$ws.Range('E50').Value = '  +2.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.75'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.36%  '
